$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-enter the formulas for B2:B40 (each references the cell below, minus a
# fixed decrement) with the new decrement (3 instead of 10.8).
$rng = $ws.Range("B2:B40")
$rng.FormulaR1C1 = "=R[1]C+3"

# B41 keeps its own (non-shared) formula referencing B42.
$ws.Range("B41").Formula = "=B42+3"

# Select the range the same way the author did when making the edit.
$sel = $ws.Range("B2:B41")
$sel.Select()
$ws.Range("B41").Activate()

$wb.Save()
